{"js": "// Move the \"Play ... Game Review\" / \"Discover ...\" pair from the end of the\n// document up to a new \"Meta description\" paragraph right under the H1\n// title, and replace the old trailing pair with a single paragraph that now\n// holds an image-generation prompt (keeping the original italic run).\n\nconst body = context.document.body;\nlet paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// --- Step 1: insert a new \"Meta description\" paragraph right after the H1 title.\nconst title = paras.items[0];\nconst metaPara = title.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nconst metaOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r/>\n            <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\">: Discover Alpha Eagle Stack N Sync: an immersive winter-themed game with unique bonuses and special symbols. Play now for free!</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nmetaPara.insertOoxml(metaOoxml, \"Replace\");\nawait context.sync();\n\n// --- Step 2: re-load paragraphs (indices shifted because of the insert above).\nparas = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst n = paras.items.length;\nconst boldTitlePara = paras.items[n - 2];  // trailing bold \"Play ... Game Review\"\nconst italicDescPara = paras.items[n - 1]; // trailing italic \"Discover ...\" paragraph\n\n// --- Step 3: drop the now-duplicated bold title paragraph entirely.\nboldTitlePara.delete();\nawait context.sync();\n\n// --- Step 4: swap the italic paragraph's text for the new image prompt,\n// keeping its run/paragraph formatting (leading empty run + italic run).\nconst descOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r/>\n            <w:r><w:rPr><w:i/></w:rPr><w:t>Create a cartoon-style feature image for &quot;Alpha Eagle Stack N Sync&quot; that features a happy Maya warrior with glasses. The image should be colorful with a winter mountain backdrop, showcasing the beautiful bald eagle as the main focus, and have the warrior holding a golden version of the S'N'S symbol with an excited expression on their face. The eagle's talons should also be featured in the image, along with other relevant symbols from the game such as feathers and salmon. The image should be eye-catching and give players an idea of the game's exciting features while also showcasing the fun and playful attitude of the game.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nitalicDescPara.insertOoxml(descOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Move the \"Play ... Game Review\" / \"Discover ...\" pair from the end of the\n# document up to a new \"Meta description\" paragraph right under the H1\n# title, and replace the old trailing pair with a single paragraph that now\n# holds an image-generation prompt (keeping the original italic run).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert a new \"Meta description\" paragraph right after the H1 title.\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n\n$metaOoxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r/>\n            <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\">: Discover Alpha Eagle Stack N Sync: an immersive winter-themed game with unique bonuses and special symbols. Play now for free!</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$metaPara.Range.InsertXML($metaOoxml, \"Replace\")\n\n# --- Step 2: locate the trailing bold \"Play ... Game Review\" paragraph and\n# the italic \"Discover ...\" paragraph that follows it (now shifted down by\n# one because of the insert above).\n$n = $d.Paragraphs.Count\n$boldTitlePara = $d.Paragraphs.Item($n - 1)\n$italicDescPara = $d.Paragraphs.Item($n)\n\n# --- Step 3: drop the now-duplicated bold title paragraph entirely.\n$boldTitlePara.Range.Delete()\n\n# --- Step 4: swap the italic paragraph's text for the new image prompt,\n# keeping its run/paragraph formatting (leading empty run + italic run).\n$n2 = $d.Paragraphs.Count\n$italicDescPara = $d.Paragraphs.Item($n2)\n\n$descOoxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r/>\n            <w:r><w:rPr><w:i/></w:rPr><w:t>Create a cartoon-style feature image for \"Alpha Eagle Stack N Sync\" that features a happy Maya warrior with glasses. The image should be colorful with a winter mountain backdrop, showcasing the beautiful bald eagle as the main focus, and have the warrior holding a golden version of the S'N'S symbol with an excited expression on their face. The eagle's talons should also be featured in the image, along with other relevant symbols from the game such as feathers and salmon. The image should be eye-catching and give players an idea of the game's exciting features while also showcasing the fun and playful attitude of the game.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$italicDescPara.Range.InsertXML($descOoxml, \"Replace\")\n"}
